$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HBAN")

# Row 12 - Gross Margin
$ws.Range("D12").Value = 1.0049
$ws.Range("E12").Value = 0.9384
$ws.Range("F12").Value = 0.8878
$ws.Range("G12").Value = 0.8253

# Row 14 - EBT margin
$ws.Range("D14").Value = 0.2057
$ws.Range("E14").Value = 0.2108
$ws.Range("F14").Value = 0.246
$ws.Range("G14").Value = 0.2934

# Row 15 - Net Profit Margin
$ws.Range("D15").Value = 0.1558
$ws.Range("E15").Value = 0.1632
$ws.Range("F15").Value = 0.1951
$ws.Range("G15").Value = 0.2364

# Row 16 - Free Cash Flow Margin
$ws.Range("B16").Value = 0.2282
$ws.Range("D16").Value = 0.4624
$ws.Range("E16").Value = 0.2919
$ws.Range("F16").Value = 0.273
$ws.Range("G16").Value = 0.2594

# Row 18 - EPS (Basic, Consolidated)
$ws.Range("B18").Value = 1.2789

# Row 19 - EPS (Basic, from Continuous Ops)
$ws.Range("B19").Value = 1.2789

# Row 23 - Operating Cash Flow Margin
$ws.Range("B23").Value = 0.2546
$ws.Range("D23").Value = 0.4851
$ws.Range("E23").Value = 0.313
$ws.Range("F23").Value = 0.2908
$ws.Range("G23").Value = 0.2783
